$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-04"

# Update the August row label (column A, row 9)
$ws.Range("A9").Value = "August (through 08-04)"

# Update August row (row 9) data values for each year column (B..I)
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = 11
$ws.Range("F9").Value = 7
$ws.Range("G9").Value = 25
$ws.Range("H9").Value = 20
$ws.Range("I9").Value = 22

# Update Total row (row 10) data values for each year column (B..I)
$ws.Range("B10").Value = 165
$ws.Range("C10").Value = 306
$ws.Range("D10").Value = 473
$ws.Range("E10").Value = 436
$ws.Range("F10").Value = 311
$ws.Range("G10").Value = 646
$ws.Range("H10").Value = 930
$ws.Range("I10").Value = 992
